$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh -- each D/E cell is forced to remain TEXT
# (leading apostrophe keeps Excel from re-parsing numeric-looking strings),
# then Style is reset to Normal so the quotePrefix formatting bit added by
# the apostrophe entry doesn't linger on the cell.

$c = $ws.Range("D2")
$c.Value = '''60.910.30'
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = '''  +0.42%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = '''2.919.22'
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = '''  +0.88%  '
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = '''  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = '''593.73'
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = '''  +1.64%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = '''145.63'
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = '''  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = '''  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = '''0.506'
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = '''  +0.69%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = '''6.86'
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = '''  +2.72%  '
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = '''  +0.63%  '
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = '''0.440'
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = '''  -1.29%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = '''0.0000226'
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = '''  +1.14%  '
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = '''33.66'
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = '''  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = '''  +0.11%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = '''3.400.38'
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = '''  +0.84%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = '''60.905.73'
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = '''  +0.57%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = '''6.71'
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = '''  -0.99%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = '''2.920.68'
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = '''  +0.90%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = '''430.27'
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = '''  +1.23%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = '''13.37'
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = '''  -1.41%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = '''0.682'
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = '''  +2.47%  '
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = '''  +0.28%  '
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = '''81.43'
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = '''  +1.86%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = '''10.99'
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = '''  +0.68%  '
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = '''  +0.56%  '
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = '''11.93'
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = '''  +0.91%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = '''1.00'
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = '''  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = '''  +5.58%  '
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = '''  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = '''  +0.43%  '
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = '''7.06'
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = '''  -1.70%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = '''26.43'
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = '''  +0.51%  '
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = '''0.108'
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = '''  +2.10%  '
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = '''0.0₃0847'
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = '''  +2.17%  '
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = '''  +0.98%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = '''5.63'
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = '''  +3.34%  '
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = '''  +1.01%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = '''1.99'
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = '''  -1.47%  '
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = '''  -1.07%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = '''0.287'
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = '''  -0.35%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = '''40.47'
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = '''  -2.45%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = '''373.65'
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = '''  +0.78%  '
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = '''  +0.42%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = '''2.712.83'
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = '''  +2.50%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = '''130.68'
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = '''  -0.96%  '
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = '''  -0.06%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = '''23.96'
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = '''  -3.38%  '
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = '''  +0.30%  '
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = '''  -2.64%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = '''0.126'
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = '''  +2.71%  '
$c.Style = "Normal"
